$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '92.046.27'
$ws.Range("E2").Value = '  +1.99%  '
$ws.Range("D3").Value = '3.165.12'
$ws.Range("E3").Value = '  +2.59%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''241.24'
$ws.Range("E5").Value = '  +3.21%  '
$ws.Range("D6").Value = '''623.64'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '''1.14'
$ws.Range("E7").Value = '  +6.01%  '
$ws.Range("D8").Value = '''0.375'
$ws.Range("E8").Value = '  +3.27%  '
$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '3.159.99'
$ws.Range("E10").Value = '  +2.36%  '
$ws.Range("D11").Value = '''0.755'
$ws.Range("E11").Value = '  +3.93%  '
$ws.Range("E12").Value = '  +4.32%  '
$ws.Range("D13").Value = '''0.0000252'
$ws.Range("E13").Value = '  -0.09%  '
$ws.Range("E14").Value = '  -2.46%  '
$ws.Range("D15").Value = '''5.53'
$ws.Range("E15").Value = '  +1.08%  '
$ws.Range("D16").Value = '91.880.81'
$ws.Range("E16").Value = '  +2.25%  '
$ws.Range("D17").Value = '3.734.17'
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("D18").Value = '3.094.95'
$ws.Range("E18").Value = '  -0.10%  '
$ws.Range("D19").Value = '''3.74'
$ws.Range("E19").Value = '  -2.40%  '
$ws.Range("D20").Value = '''15.76'
$ws.Range("E20").Value = '  +12.55%  '
$ws.Range("D21").Value = '''0.0000212'
$ws.Range("E21").Value = '  -1.74%  '
$ws.Range("D22").Value = '''5.82'
$ws.Range("E22").Value = '  +4.80%  '
$ws.Range("D23").Value = '''446.78'
$ws.Range("E23").Value = '  +2.66%  '
$ws.Range("D24").Value = '''9.34'
$ws.Range("E24").Value = '  +4.97%  '
$ws.Range("D25").Value = '''6.03'
$ws.Range("E25").Value = '  +4.97%  '
$ws.Range("D26").Value = '''90.26'
$ws.Range("E26").Value = '  +2.10%  '
$ws.Range("E27").Value = '  -0.95%  '
$ws.Range("D28").Value = '3.259.02'
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '''0.139'
$ws.Range("E30").Value = '  +53.40%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '''0.242'
$ws.Range("E31").Value = '  +22.06%  '
$ws.Range("B32").Value = 'Cronos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D32").Value = '''0.173'
$ws.Range("E32").Value = '  +9.00%  '
$ws.Range("D33").Value = '''9.34'
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("D34").Value = '''0.170'
$ws.Range("E34").Value = '  +10.10%  '
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").Value = '''0.988'
$ws.Range("E35").Value = '  -0.88%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D36").Value = '''8.06'
$ws.Range("E36").Value = '  +12.80%  '
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").Value = '''26.79'
$ws.Range("E37").Value = '  +3.89%  '
$ws.Range("B38").Value = 'MantraDAO'
$ws.Range("C38").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D38").Value = '''4.27'
$ws.Range("E38").Value = '  +24.37%  '
$ws.Range("D39").Value = '''510.14'
$ws.Range("E39").Value = '  +1.15%  '
$ws.Range("D40").Value = '''1.95'
$ws.Range("E40").Value = '  +2.22%  '
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").Value = '''1.32'
$ws.Range("E41").Value = '  +2.33%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '''3.52'
$ws.Range("E42").Value = '  -8.32%  '
$ws.Range("D43").Value = '''0.439'
$ws.Range("E43").Value = '  +8.17%  '
$ws.Range("D44").Value = '''22.21'
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '''0.715'
$ws.Range("E46").Value = '  +3.64%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = '''1.95'
$ws.Range("E47").Value = '  +2.03%  '
$ws.Range("D48").Value = '''154.23'
$ws.Range("E48").Value = '  +1.60%  '
$ws.Range("D49").Value = '''1.37'
$ws.Range("E49").Value = '  +2.30%  '
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").Value = '''4.48'
$ws.Range("E50").Value = '  +1.48%  '
$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").Value = '''44.44'
$ws.Range("E51").Value = '  -1.00%  '
